# Scheduled runner update: refresh market price / profit figures across Sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5 (Leve Item ID 5503)
$ws.Range("H5").Value = 63
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
# Row 29 (Leve Item ID 4575)
$ws.Range("H29").Value = 76.5
$ws.Range("I29").Value = 76.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 229.5
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = 51.5
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 4274
$ws.Range("J137").Value = 5400
$ws.Range("L137").Value = 16200
$ws.Range("N137").Value = -21300
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2878.9092
$ws.Range("J138").Value = 4101.9414
$ws.Range("L138").Value = 12305.8242
$ws.Range("N138").Value = -22585.8242
# Row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 8991.429
$ws.Range("I141").Value = 9030.556
$ws.Range("K141").Value = 27091.668
$ws.Range("M141").Value = -21911.668

$ws = $wb.Worksheets.Item("ARM")
# Row 31 (Leve Item ID 19533)
$ws.Range("H31").Value = 12143.5
$ws.Range("I31").Value = 11382.091
$ws.Range("K31").Value = 11382.091
$ws.Range("M31").Value = -11088.091
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 3599
$ws.Range("I45").Value = 3599
$ws.Range("K45").Value = 3599
$ws.Range("M45").Value = -3222
# Row 51 (Leve Item ID 3858)
$ws.Range("H51").Value = 52500
$ws.Range("J51").Value = 52500
$ws.Range("L51").Value = 52500
$ws.Range("N51").Value = -54012
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 3034.5386
$ws.Range("I61").Value = 1911.875
$ws.Range("K61").Value = 1911.875
$ws.Range("M61").Value = -1699.875
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 3034.5386
$ws.Range("I136").Value = 1911.875
$ws.Range("K136").Value = 5735.625
$ws.Range("M136").Value = -3185.625

$ws = $wb.Worksheets.Item("BSM")
# Row 70 (Leve Item ID 15553)
$ws.Range("H70").Value = 80000
$ws.Range("J70").Value = 80000
$ws.Range("L70").Value = 80000
$ws.Range("N70").Value = -80586
# Row 73 (Leve Item ID 15553)
$ws.Range("H73").Value = 80000
$ws.Range("J73").Value = 80000
$ws.Range("L73").Value = 80000
$ws.Range("N73").Value = -82028
# Row 102 (Leve Item ID 19565)
$ws.Range("H102").Value = 30361.092
$ws.Range("I102").Value = 14562.777
$ws.Range("J102").Value = 101453.5
$ws.Range("K102").Value = 14562.777
$ws.Range("L102").Value = 101453.5
$ws.Range("M102").Value = -11317.777
$ws.Range("N102").Value = -107943.5

$ws = $wb.Worksheets.Item("CRP")
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 1996.5238
$ws.Range("I134").Value = 1223
$ws.Range("K134").Value = 3669
$ws.Range("M134").Value = -1134
# Row 140 (Leve Item ID 42455)
$ws.Range("H140").Value = 46142.785
$ws.Range("J140").Value = 46142.785
$ws.Range("L140").Value = 46142.785
$ws.Range("N140").Value = -56502.785

$ws = $wb.Worksheets.Item("CUL")
# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 4249.75
$ws.Range("J68").Value = 4249.75
$ws.Range("L68").Value = 12749.25
$ws.Range("N68").Value = -14371.25
# Row 69 (Leve Item ID 12850)
$ws.Range("H69").Value = 111115090
$ws.Range("J69").Value = 111115090
$ws.Range("L69").Value = 333345270
$ws.Range("N69").Value = -333346892
# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 4249.75
$ws.Range("J71").Value = 4249.75
$ws.Range("L71").Value = 38247.75
$ws.Range("N71").Value = -46359.75
# Row 72 (Leve Item ID 12850)
$ws.Range("H72").Value = 111115090
$ws.Range("J72").Value = 111115090
$ws.Range("L72").Value = 1000035810
$ws.Range("N72").Value = -1000043922
# Row 80 (Leve Item ID 12890)
$ws.Range("H80").Value = 7772.727
$ws.Range("J80").Value = 8083.3335
$ws.Range("L80").Value = 24250.0005
$ws.Range("N80").Value = -26122.0005
# Row 83 (Leve Item ID 12890)
$ws.Range("H83").Value = 7772.727
$ws.Range("J83").Value = 8083.3335
$ws.Range("L83").Value = 72750.0015
$ws.Range("N83").Value = -82110.0015
# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 1959.3024
$ws.Range("J107").Value = 1843.909
$ws.Range("L107").Value = 5531.727000000001
$ws.Range("N107").Value = -9371.727000000001
# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 7141.0835
$ws.Range("J122").Value = 10481.875
$ws.Range("L122").Value = 94336.875
$ws.Range("N122").Value = -99236.875
# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 2021.5
$ws.Range("I132").Value = 1978.4
$ws.Range("K132").Value = 17805.6
$ws.Range("M132").Value = -15275.6
# Row 136 (Leve Item ID 44093)
$ws.Range("H136").Value = 2471.75
$ws.Range("I136").Value = 949.7778
$ws.Range("K136").Value = 2849.3334
$ws.Range("M136").Value = 2250.6666
# Row 137 (Leve Item ID 44088)
$ws.Range("H137").Value = 2146.2222
$ws.Range("J137").Value = 4281
$ws.Range("L137").Value = 12843
$ws.Range("N137").Value = -23043

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 5993.3
$ws.Range("I122").Value = 5655.6665
$ws.Range("J122").Value = 6499.75
$ws.Range("K122").Value = 16966.9995
$ws.Range("L122").Value = 19499.25
$ws.Range("M122").Value = -14516.9995
$ws.Range("N122").Value = -24399.25
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 7762.6577
$ws.Range("I132").Value = 7054.125
$ws.Range("K132").Value = 21162.375
$ws.Range("M132").Value = -18632.375
# Row 137 (Leve Item ID 43226)
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0
# Row 139 (Leve Item ID 42373)
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 2500.3333
$ws.Range("I7").Value = 2435.6667
$ws.Range("J7").Value = 2565
$ws.Range("K7").Value = 2435.6667
$ws.Range("L7").Value = 2565
$ws.Range("M7").Value = -2323.6667
$ws.Range("N7").Value = -2789
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 25000202
$ws.Range("I16").Value = 25000202
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 25000202
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -25000032
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 6993654.5
$ws.Range("I22").Value = 15152080
$ws.Range("J22").Value = 718.7143
$ws.Range("K22").Value = 15152080
$ws.Range("L22").Value = 718.7143
$ws.Range("M22").Value = -15151785
$ws.Range("N22").Value = -1308.7143
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 6993654.5
$ws.Range("I27").Value = 15152080
$ws.Range("J27").Value = 718.7143
$ws.Range("K27").Value = 15152080
$ws.Range("L27").Value = 718.7143
$ws.Range("M27").Value = -15151973
$ws.Range("N27").Value = -932.7143
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 419.24
$ws.Range("I55").Value = 495.2353
$ws.Range("J55").Value = 257.75
$ws.Range("K55").Value = 495.2353
$ws.Range("L55").Value = 257.75
$ws.Range("M55").Value = -322.2353
$ws.Range("N55").Value = -603.75
# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 1397.1111
$ws.Range("I61").Value = 1434.25
$ws.Range("K61").Value = 1434.25
$ws.Range("M61").Value = -1232.25
# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 2398
$ws.Range("I93").Value = 2296.5
$ws.Range("J93").Value = 2499.5
$ws.Range("K93").Value = 2296.5
$ws.Range("L93").Value = 2499.5
$ws.Range("M93").Value = -1048.5
$ws.Range("N93").Value = -4995.5
# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 1397.1111
$ws.Range("I113").Value = 1434.25
$ws.Range("K113").Value = 1434.25
$ws.Range("M113").Value = 735.75
# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 2500.3333
$ws.Range("I126").Value = 2435.6667
$ws.Range("J126").Value = 2565
$ws.Range("K126").Value = 7307.000100000001
$ws.Range("L126").Value = 7695
$ws.Range("M126").Value = -4837.000100000001
$ws.Range("N126").Value = -12635

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1781.9117
$ws.Range("I132").Value = 1635.1936
$ws.Range("K132").Value = 4905.5808
$ws.Range("M132").Value = -2375.5808
